# ---------------------------------------------------------------------------
# Refined metadata to be additional tab
#
# 1) Refresh the "panel_query_time" timestamps recorded in the `data` sheet
#    (column F) for every gene row.
# 2) Add a new "metadata" worksheet (placed right after "data") describing
#    the panelapp query the data was pulled from.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) refresh data!F2:F59 panel_query_time values -------------------------
$timestamps = @{
    2  = "2021-10-05 14:33:33.496362"
    3  = "2021-10-05 14:33:33.496370"
    4  = "2021-10-05 14:33:33.496373"
    5  = "2021-10-05 14:33:33.496376"
    6  = "2021-10-05 14:33:33.496379"
    7  = "2021-10-05 14:33:33.496381"
    8  = "2021-10-05 14:33:33.496384"
    9  = "2021-10-05 14:33:33.496386"
    10 = "2021-10-05 14:33:33.496389"
    11 = "2021-10-05 14:33:33.496392"
    12 = "2021-10-05 14:33:33.496394"
    13 = "2021-10-05 14:33:33.496397"
    14 = "2021-10-05 14:33:33.496399"
    15 = "2021-10-05 14:33:33.496402"
    16 = "2021-10-05 14:33:33.496404"
    17 = "2021-10-05 14:33:33.496406"
    18 = "2021-10-05 14:33:33.496409"
    19 = "2021-10-05 14:33:33.496412"
    20 = "2021-10-05 14:33:33.496415"
    21 = "2021-10-05 14:33:33.496417"
    22 = "2021-10-05 14:33:33.496420"
    23 = "2021-10-05 14:33:33.496422"
    24 = "2021-10-05 14:33:33.496424"
    25 = "2021-10-05 14:33:33.496427"
    26 = "2021-10-05 14:33:33.496430"
    27 = "2021-10-05 14:33:33.496432"
    28 = "2021-10-05 14:33:33.496435"
    29 = "2021-10-05 14:33:33.496437"
    30 = "2021-10-05 14:33:33.496440"
    31 = "2021-10-05 14:33:33.496442"
    32 = "2021-10-05 14:33:33.496445"
    33 = "2021-10-05 14:33:33.496447"
    34 = "2021-10-05 14:33:33.496450"
    35 = "2021-10-05 14:33:33.496453"
    36 = "2021-10-05 14:33:33.496455"
    37 = "2021-10-05 14:33:33.496458"
    38 = "2021-10-05 14:33:33.496460"
    39 = "2021-10-05 14:33:33.496463"
    40 = "2021-10-05 14:33:33.496465"
    41 = "2021-10-05 14:33:33.496468"
    42 = "2021-10-05 14:33:33.496471"
    43 = "2021-10-05 14:33:33.496473"
    44 = "2021-10-05 14:33:33.496476"
    45 = "2021-10-05 14:33:33.496478"
    46 = "2021-10-05 14:33:33.496481"
    47 = "2021-10-05 14:33:33.496483"
    48 = "2021-10-05 14:33:33.496486"
    49 = "2021-10-05 14:33:33.496488"
    50 = "2021-10-05 14:33:33.496491"
    51 = "2021-10-05 14:33:33.496493"
    52 = "2021-10-05 14:33:33.496496"
    53 = "2021-10-05 14:33:33.496499"
    54 = "2021-10-05 14:33:33.496502"
    55 = "2021-10-05 14:33:33.496504"
    56 = "2021-10-05 14:33:33.496507"
    57 = "2021-10-05 14:33:33.496509"
    58 = "2021-10-05 14:33:33.496512"
    59 = "2021-10-05 14:33:33.496515"
}

foreach ($row in $timestamps.Keys) {
    $ws1.Range("F$row").Value = $timestamps[$row]
}

# --- 2) add the "metadata" worksheet right after "data" ---------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "metadata"

# Header row (bold, bordered, centered -- matches the "data" sheet header style)
$headerRange = $newSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row describing the panel query
$a2 = $newSheet.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

$newSheet.Range("B2").Value = "Congenital anomalies of the kidney and urinary tract (CAKUT) Syndromic"
$newSheet.Range("C2").Value = 63
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.89"
$newSheet.Range("E2").Value = "2021-10-04T06:48:28.340886Z"
$newSheet.Range("F2").Value = "2021-10-05 14:33:33.492936"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/63/?format=json"

# Keep "data" as the active/selected sheet, as in the original workbook.
$ws1.Activate()

Write-Output "sheets=$($wb.Worksheets.Count) sheet1=$($wb.Worksheets.Item(1).Name) sheet2=$($wb.Worksheets.Item(2).Name) active=$($wb.ActiveSheet.Name)"
